$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Periodo Mora" column (E) for rows 16-22 is re-ordered: the list of periods
# 1902-1908 now runs ascending top-to-bottom instead of descending.
$ws.Range("E16").Value = "1902"
$ws.Range("E17").Value = "1903"
$ws.Range("E18").Value = "1904"
$ws.Range("E19").Value = "1905"
$ws.Range("E20").Value = "1906"
$ws.Range("E21").Value = "1907"
$ws.Range("E22").Value = "1908"

# "Valor Mora" column (F) updated values per row.
$ws.Range("F16").Value = 33125
$ws.Range("F17").Value = 33125
$ws.Range("F18").Value = 31249
$ws.Range("F19").Value = 31249
$ws.Range("F20").Value = 31249
$ws.Range("F21").Value = 31249
$ws.Range("F22").Value = 19791

# "Salario Basico" column (G) updated uniformly for every row.
$ws.Range("G16:G22").Value = 781242
